# Pedido 691b59eb90ee710b08f4ffcc — add new order row (row 14) to "Productos"
# and drop the stale empty placeholder cells (F13/G13/L13) from row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: Optimizador / Unidades Optimizador / Cargador VE were blank
#     placeholders left over from a prior export; clear them out entirely.
$ws.Range("F13").ClearContents()
$ws.Range("G13").ClearContents()
$ws.Range("L13").ClearContents()

# --- Row 14: new order for "Juan Carlos 55"
$row = 14

$ws.Cells.Item($row, 1).Value = 2663
$ws.Cells.Item($row, 2).Value = "Juan Carlos 55"
$ws.Cells.Item($row, 3).Value = "Estructura coplanar NOVOTEGRA"
$ws.Cells.Item($row, 4).Value = "MODULO FV JA SOLAR 535WP BLACK FRAME BIFACIAL 120 CELDAS"

# Unidades Estructura/Paneles is stored as text ("12"), not a number -
# force text formatting before assigning so Excel doesn't coerce it.
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "12"

# Optimizador / Unidades Optimizador - no optimizer on this order.
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = ""

$ws.Cells.Item($row, 8).Value = "GOODWE GW6000-ES-20 híbrido monofásico"

$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = "1"

$ws.Cells.Item($row, 10).Value = "GOODWE Batería Lynx Home U G3 5,12 kWh"

$ws.Cells.Item($row, 11).NumberFormat = "@"
$ws.Cells.Item($row, 11).Value = "1"

# Cargador VE - none on this order.
$ws.Cells.Item($row, 12).Value = ""

$ws.Cells.Item($row, 13).Value = "Sí"
$ws.Cells.Item($row, 14).Value = "2025-06-25T13:59:26.816Z"
